# Weekly data refresh: a new daily-price record is inserted at row 24,
# pushing the existing historical rows (old rows 24-123) down by one row
# (new rows 25-124). This mirrors Excel's normal "insert row" behavior
# which shifts down & carries formatting from the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 24 (everything below shifts down one row).
$ws.Rows("24:24").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A24").Value = 8
$ws.Range("B24").Value = "Terminal La Palmera de La Serena"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = (Get-Date -Year 2023 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 100114007
$ws.Range("G24").Value = "Jengibre"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 440
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16500
$ws.Range("N24").Value = "$/caja 13 kilos"
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 1269
$ws.Range("Q24").Value = 13
$ws.Range("R24").Value = "Hortaliza"
